$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.013.29"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.899.52"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.79%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.08"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5011"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3924"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.69%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.09401"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  -0.37%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "41.93"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.49%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.360"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.903.04"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.32%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.304"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001116"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "92.43"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06583"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("E21").Value = "  +0.01%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.196"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.047.31"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.29"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.310"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.383"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.625"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.119.53"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "20.85"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "156.94"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "126.62"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.078"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.1062"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.608"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.611"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "9.635"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06610"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02413"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2175"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.226"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.244"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.991"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6326"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "11.37"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "13.26"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5958"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.713"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("B49").Value = "WEMIXTOKEN"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.274"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.015"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "123.19"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
